$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 84

$ws.Cells.Item($newRow, 1).Value = "Partou"
$ws.Cells.Item($newRow, 2).Value = "Partou Hetty Blokweg 2"
$ws.Cells.Item($newRow, 3).Value = "KDV"
$ws.Cells.Item($newRow, 4).Value = "'2024-06-12"
$ws.Cells.Item($newRow, 5).Value = 0
$ws.Cells.Item($newRow, 6).Value = 0
$ws.Cells.Item($newRow, 7).Value = 0
$ws.Cells.Item($newRow, 8).Value = 0
$ws.Cells.Item($newRow, 9).Value = 0
$ws.Cells.Item($newRow, 10).Value = 0
